# Normalize every product's category (column I, rows 2-21) to "sed".
# This matches the source data re-import/refactor described in the commit
# message: the category values that used to vary per-row ("sit", "rerum",
# "nobis", "ut", "quae", "aut", "quo", "sed", ...) are all consolidated to
# a single value, "sed". Excel automatically drops the now-unreferenced
# shared strings ("sit", "rerum", "nobis", "ut", "quae") from the shared
# string table on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2:I21").Value = "sed"

# Reflect the selection left behind by the edit (selecting the column that
# was just changed before saving).
$ws.Range("I2:I21").Select()
